# Weekly update: a new week's price record is inserted at the top of the
# data block (row 15), pushing all existing data rows (15-23) down by one
# (to 16-24). This mirrors the "Fruta / hortaliza, semanal" commit, which
# adds the latest week's observation while keeping older weeks below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 15-23 down to 16-24, inserting a blank row 15.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new week's record.
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 44839
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 100112036
$ws.Range("G15").Value = "Caigua"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 16000
$ws.Range("L15").Value = 16000
$ws.Range("M15").Value = 16000
$ws.Range("N15").Value = "$/caja 15 kilos"
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 1067
$ws.Range("Q15").Value = 15
$ws.Range("R15").Value = "Hortaliza"
